$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange / UpDown for row 31
$ws.Range("X31").Value = -0.21000000000000085
$ws.Range("Y31").Value = "Down"

# New row 32 data
$ws.Range("A32").Value = 42651.425173611111
$ws.Range("A32").NumberFormat = "m/d/yy h:mm"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "Neutral"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = "Random"
$ws.Range("Q32").Value = 47.963765586266284
$ws.Range("R32").Value = 0.49
$ws.Range("S32").Value = 0.0521
$ws.Range("S32").NumberFormat = "0.00%"
$ws.Range("T32").Value = -0.0214
$ws.Range("T32").NumberFormat = "0.00%"
$ws.Range("U32").Value = 2.2799999999999998
$ws.Range("V32").Value = "N/A"
$ws.Range("W32").Value = 0
